$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.8715449999999999
$ws.Range("H2").Value = 2.614635
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 44.57468866666667
$ws.Range("N2").Value = 133.724066
$ws.Range("O2").Value = 0.5438239100642482
$ws.Range("P2").Value = 0.5438239100642482
$ws.Range("Q2").Value = 38.84884703399
$ws.Range("R2").Value = 349.6396233059099
$ws.Range("S2").Value = 0.5438239100642482
$ws.Range("T2").Value = 0.5438239100642482

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.8715449999999999
$ws.Range("H3").Value = 2.614635
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 27.876023
$ws.Range("N3").Value = 83.628069
$ws.Range("O3").Value = 0.3400954281086753
$ws.Range("P3").Value = 0.3400954281086753
$ws.Range("Q3").Value = 24.295208465535
$ws.Range("R3").Value = 218.656876189815
$ws.Range("S3").Value = 0.3400954281086753
$ws.Range("T3").Value = 0.3400954281086753

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.8715449999999999
$ws.Range("H4").Value = 2.614635
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 9.514586
$ws.Range("N4").Value = 28.543758
$ws.Range("O4").Value = 0.1160806618270766
$ws.Range("P4").Value = 0.1160806618270765
$ws.Range("Q4").Value = 8.292389855369999
$ws.Range("R4").Value = 74.63150869832998
$ws.Range("S4").Value = 0.1160806618270766
$ws.Range("T4").Value = 0.1160806618270765
